$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (mirrors "Insert Sheet"
# placed at the end of the tab strip), then rename it to Sheet2.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Sheet2"

# Row 1 - pasted status line from the training log website.
$ws.Range("A1").Value = "今天7:31 2011-04-01查我的学车资料状态已经是已结业了"

# Row 4 - summary fields scraped from the page (A-F, then J-L).
$ws.Range("A4").Value = "[查看]"
$ws.Range("B4").Value = " 张昀"
$ws.Range("C4").Value = " C1"
$ws.Range("D4").Value = " 已结业"
$ws.Range("E4").Value = " 初次申请"
$ws.Range("F4").Value = " 2010-10-25 00:00:00"
$ws.Range("J4").Value = " 广州市穗通驾驶员培训有限公司"
$ws.Range("K4").Value = " JP440106000014"
$ws.Range("L4").Value = " S10111016"

# Rows 5-11 - stage / hours breakdown.
$ws.Range("A5").Value = "当前处于阶段"
$ws.Range("A6").Value = "大纲要求培训时长为:"
$ws.Range("A7").Value = "理论"
$ws.Range("A8").Value = "实操"
$ws.Range("A9").Value = "实际完成："
$ws.Range("A10").Value = "理论"
$ws.Range("A11").Value = "实操"

# Row 12 - table header.
$ws.Range("A12").Value = "学员名称"
$ws.Range("B12").Value = "培训车型"
$ws.Range("C12").Value = "培训状态"
$ws.Range("D12").Value = "培训类型"
$ws.Range("E12").Value = "入学时间"

# Row 13 - table data.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Select the pasted block and make this new sheet the active tab, matching
# the saved view state (tabSelected moves from Sheet1 to Sheet2).
$ws.Range("A1:L13").Select()
$ws.Activate()
